$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should match the styling of the
# existing header cells (bold font, thin border, centered/top aligned).
# Copy the format from H1 (an existing styled header cell) into I1:J1, then
# set the header text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 (col I) and IF (col J) columns, rows 2-57.
$iVals = @(4,5,7,9,6,9,9,7,5,6,6,7,9,7,7,7,8,6,7,6,6,7,10,9,8,7,6,7,7,7,7,6,8,7,9,8,8,8,6,7,9,8,5,7,6,6,5,7,2,1,3,7,3,5,6,3)
$jVals = @(7,7,8,9,6,9,9,7,6,7,6,7,9,7,7,7,8,6,7,6,6,8,10,9,8,7,7,7,7,8,7,7,9,7,9,8,8,8,6,8,10,8,7,8,7,6,7,9,3,3,3,7,3,5,6,3)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$i]
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
}
